$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.738.40'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.850.66'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -2.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.013'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4330'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3773'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07387'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8839'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').Value = '1.863.41'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.754'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07137'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.016'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009030'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').Value = '27.770.60'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.268'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').Value = '2.086.69'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.031'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.141'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.416'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08970'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.237'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7777'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.577'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.926'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.147'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.014'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05334'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.161'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.861'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5193'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1685'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.944'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '110.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.717'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4744'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06511'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.014'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.897'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
